# Update the EPEX Spot prices workbook with the latest day of data.
#
# Sheet "Prix Spot": add a new date column (BR) "22-aug" with 24 hourly prices.
# Sheet "Gaz":       add a new row (67) for date 2025-08-20.
# Sheet "CO2":       add a new row (67) for date 2025-08-20.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Prix Spot" -- new column BR ("22-aug")
# ---------------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the header style/formatting from the previous date column (BQ1) onto
# the new one (BR1), then set its text.
$wsSpot.Range("BQ1").Copy($wsSpot.Range("BR1"))
$wsSpot.Range("BR1").Value = "22-aug"

$spotValues = @(
    66.38,
    51.66,
    48.66,
    41.88,
    31.51,
    36.41,
    71.03,
    70.88,
    75.09999999999999,
    49.67,
    15.57,
    5.11,
    19.02,
    6.83,
    5.37,
    10.07,
    23.6,
    26.45,
    73.48,
    84.56,
    104.53,
    103.49,
    101.59,
    94.37
)

for ($i = 0; $i -lt $spotValues.Count; $i++) {
    $row = $i + 2
    $wsSpot.Cells.Item($row, 70).Value = $spotValues[$i]
}

# ---------------------------------------------------------------------------
# Sheet 2: "Gaz" -- new row 67 (2025-08-20, 30.85)
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the date cell to be stored as plain text (matching the rest of the
# column) instead of being auto-converted into a date serial number.
$wsGaz.Range("A67").NumberFormat = "@"
$wsGaz.Range("A67").Value = "2025-08-20"
$wsGaz.Range("A67").Style = "Normal"

$wsGaz.Range("B67").Value = 30.85

# ---------------------------------------------------------------------------
# Sheet 3: "CO2" -- new row 67 (2025-08-20, 70.3)
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A67").NumberFormat = "@"
$wsCo2.Range("A67").Value = "2025-08-20"
$wsCo2.Range("A67").Style = "Normal"

$wsCo2.Range("B67").Value = 70.3
